# Actualización automática 2025-07-14 13:30:09
$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasGrupo.Range("M18").Value = 350.46
$wsVentasGrupo.Range("O18").Value = 62.21

# --- Sheet: VENTA MENSUAL ---
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F18").Value = 514.49
$wsVentaMensual.Range("F32").Value = 2979.77

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 15 - PORCELANATO
$wsCumplimiento.Range("D15").Value = 479.46
$wsCumplimiento.Range("E15").Value = 22979.36
$wsCumplimiento.Range("F15").Value = 0.02043836817026602

# Row 17 - SAL SOLUBLE
$wsCumplimiento.Range("D17").Value = 60.44
$wsCumplimiento.Range("E17").Value = 1539.56
$wsCumplimiento.Range("F17").Value = 0.037775

# Row 18 - TOTAL
$wsCumplimiento.Range("D18").Value = 2969.49
$wsCumplimiento.Range("E18").Value = 30965.22607548726
$wsCumplimiento.Range("F18").Value = 0.08750596272544066
